$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29500
$ws.Range("J21").Value = 29500
$ws.Range("L21").Value = 29500
$ws.Range("N21").Value = -30436
$ws.Range("H23").Value = 29500
$ws.Range("J23").Value = 29500
$ws.Range("L23").Value = 29500
$ws.Range("N23").Value = -29968
$ws.Range("H40").Value = 5923.2
$ws.Range("J40").Value = 7233
$ws.Range("L40").Value = 7233
$ws.Range("N40").Value = -7583
$ws.Range("H86").Value = 666667460
$ws.Range("J86").Value = 2399
$ws.Range("L86").Value = 2399
$ws.Range("N86").Value = -4645
$ws.Range("H89").Value = 666667460
$ws.Range("J89").Value = 2399
$ws.Range("L89").Value = 11995
$ws.Range("N89").Value = -23227
$ws.Range("H112").Value = 2453.7778
$ws.Range("J112").Value = 2492.2354
$ws.Range("L112").Value = 7476.706200000001
$ws.Range("N112").Value = -9692.706200000001
$ws.Range("H132").Value = 6087.9375
$ws.Range("I132").Value = 6635
$ws.Range("K132").Value = 19905
$ws.Range("M132").Value = -17375
$ws.Range("H138").Value = 362281.44
$ws.Range("J138").Value = 415135.7
$ws.Range("L138").Value = 1245407.1
$ws.Range("N138").Value = -1255687.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1560
$ws.Range("I2").Value = 1416.44
$ws.Range("J2").Value = 2158.1667
$ws.Range("K2").Value = 1416.44
$ws.Range("L2").Value = 2158.1667
$ws.Range("M2").Value = -1303.44
$ws.Range("N2").Value = -2384.1667
$ws.Range("H32").Value = 10705.033
$ws.Range("I32").Value = 7327.075
$ws.Range("J32").Value = 35272
$ws.Range("K32").Value = 7327.075
$ws.Range("L32").Value = 35272
$ws.Range("M32").Value = -7040.075
$ws.Range("N32").Value = -35846
$ws.Range("H45").Value = 14876.567
$ws.Range("I45").Value = 21092.61
$ws.Range("J45").Value = 4664.5
$ws.Range("K45").Value = 21092.61
$ws.Range("L45").Value = 4664.5
$ws.Range("M45").Value = -20715.61
$ws.Range("N45").Value = -5418.5
$ws.Range("H46").Value = 1430.3334
$ws.Range("I46").Value = 1396.5
$ws.Range("J46").Value = 1498
$ws.Range("K46").Value = 1396.5
$ws.Range("L46").Value = 1498
$ws.Range("M46").Value = -1077.5
$ws.Range("N46").Value = -2136
$ws.Range("H61").Value = 5836.5
$ws.Range("I61").Value = 6826
$ws.Range("K61").Value = 6826
$ws.Range("M61").Value = -6614
$ws.Range("H116").Value = 1560
$ws.Range("I116").Value = 1416.44
$ws.Range("J116").Value = 2158.1667
$ws.Range("K116").Value = 1416.44
$ws.Range("L116").Value = 2158.1667
$ws.Range("M116").Value = 877.5599999999999
$ws.Range("N116").Value = -6746.1667
$ws.Range("H122").Value = 8144.1113
$ws.Range("J122").Value = 8144.1113
$ws.Range("L122").Value = 24432.3339
$ws.Range("N122").Value = -29332.3339
$ws.Range("H132").Value = 2659.6572
$ws.Range("I132").Value = 2011.4286
$ws.Range("J132").Value = 3632
$ws.Range("K132").Value = 6034.2858
$ws.Range("L132").Value = 10896
$ws.Range("M132").Value = -3504.2858
$ws.Range("N132").Value = -15956
$ws.Range("H134").Value = 89694
$ws.Range("J134").Value = 89694
$ws.Range("L134").Value = 89694
$ws.Range("N134").Value = -99834
$ws.Range("H136").Value = 5836.5
$ws.Range("I136").Value = 6826
$ws.Range("K136").Value = 20478
$ws.Range("M136").Value = -17928
$ws.Range("H138").Value = 102819.664
$ws.Range("J138").Value = 102819.664
$ws.Range("L138").Value = 102819.664
$ws.Range("N138").Value = -113099.664
$ws.Range("H140").Value = 88711.86
$ws.Range("J140").Value = 88711.86
$ws.Range("L140").Value = 88711.86
$ws.Range("N140").Value = -99071.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1560
$ws.Range("I3").Value = 1416.44
$ws.Range("J3").Value = 2158.1667
$ws.Range("K3").Value = 1416.44
$ws.Range("L3").Value = 2158.1667
$ws.Range("M3").Value = -1302.44
$ws.Range("N3").Value = -2386.1667
$ws.Range("H94").Value = 111111830
$ws.Range("I94").Value = 142857780
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 142857780
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = -142857329
$ws.Range("N94").Value = -1901
$ws.Range("H99").Value = 8814.5
$ws.Range("I99").Value = 8139.4165
$ws.Range("J99").Value = 9827.125
$ws.Range("K99").Value = 8139.4165
$ws.Range("L99").Value = 9827.125
$ws.Range("M99").Value = -6641.4165
$ws.Range("N99").Value = -12823.125
$ws.Range("H134").Value = 3163.3572
$ws.Range("I134").Value = 2843.0952
$ws.Range("K134").Value = 8529.285600000001
$ws.Range("M134").Value = -5994.285600000001
$ws.Range("H135").Value = 116635
$ws.Range("J135").Value = 116635
$ws.Range("L135").Value = 116635
$ws.Range("N135").Value = -126775
$ws.Range("H139").Value = 104989.664
$ws.Range("J139").Value = 104989.664
$ws.Range("L139").Value = 104989.664
$ws.Range("N139").Value = -115269.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8120.1743
$ws.Range("I31").Value = 5725.7734
$ws.Range("K31").Value = 5725.7734
$ws.Range("M31").Value = -5430.7734
$ws.Range("H34").Value = 8120.1743
$ws.Range("I34").Value = 5725.7734
$ws.Range("K34").Value = 5725.7734
$ws.Range("M34").Value = -5523.7734
$ws.Range("H58").Value = 3187.9512
$ws.Range("J58").Value = 4749.3076
$ws.Range("L58").Value = 4749.3076
$ws.Range("N58").Value = -5155.3076
$ws.Range("H132").Value = 11908905
$ws.Range("I132").Value = 13516784
$ws.Range("K132").Value = 40550352
$ws.Range("M132").Value = -40547822
$ws.Range("H134").Value = 3079.4722
$ws.Range("I134").Value = 2186.1785
$ws.Range("J134").Value = 6206
$ws.Range("K134").Value = 6558.5355
$ws.Range("L134").Value = 18618
$ws.Range("M134").Value = -4023.5355
$ws.Range("N134").Value = -23688
$ws.Range("H136").Value = 3187.9512
$ws.Range("J136").Value = 4749.3076
$ws.Range("L136").Value = 14247.9228
$ws.Range("N136").Value = -19347.9228
$ws.Range("H141").Value = 481496.78
$ws.Range("J141").Value = 481496.78
$ws.Range("L141").Value = 481496.78
$ws.Range("N141").Value = -491856.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1962.4
$ws.Range("I70").Value = 828
$ws.Range("K70").Value = 2484
$ws.Range("M70").Value = -2169
$ws.Range("H73").Value = 1962.4
$ws.Range("I73").Value = 828
$ws.Range("K73").Value = 2484
$ws.Range("M73").Value = -1392
$ws.Range("H122").Value = 1817.8077
$ws.Range("I122").Value = 1799
$ws.Range("J122").Value = 1823.45
$ws.Range("K122").Value = 16191
$ws.Range("L122").Value = 16411.05
$ws.Range("M122").Value = -13741
$ws.Range("N122").Value = -21311.05
$ws.Range("H129").Value = 6346.3335
$ws.Range("J129").Value = 3852.5557
$ws.Range("L129").Value = 11557.6671
$ws.Range("N129").Value = -21557.6671
$ws.Range("H138").Value = 4507.375
$ws.Range("I138").Value = 4507.375
$ws.Range("K138").Value = 13522.125
$ws.Range("M138").Value = -8382.125
$ws.Range("H141").Value = 35039.645
$ws.Range("I141").Value = 10277.5
$ws.Range("J141").Value = 39166.668
$ws.Range("K141").Value = 30832.5
$ws.Range("L141").Value = 117500.004
$ws.Range("M141").Value = -25652.5
$ws.Range("N141").Value = -127860.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 249994.5
$ws.Range("J133").Value = 249994.5
$ws.Range("L133").Value = 249994.5
$ws.Range("N133").Value = -260114.5
$ws.Range("H138").Value = 99725
$ws.Range("J138").Value = 99725
$ws.Range("L138").Value = 99725
$ws.Range("N138").Value = -110005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 864.9091
$ws.Range("I16").Value = 851.45
$ws.Range("K16").Value = 851.45
$ws.Range("M16").Value = -681.45
$ws.Range("H46").Value = 1695.8334
$ws.Range("I46").Value = 1100
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 1750
$ws.Range("M46").Value = -912
$ws.Range("N46").Value = -2126
$ws.Range("H55").Value = 699.6667
$ws.Range("I55").Value = 699.6667
$ws.Range("K55").Value = 699.6667
$ws.Range("M55").Value = -526.6667
$ws.Range("H122").Value = 6845.923
$ws.Range("I122").Value = 6139.1333
$ws.Range("J122").Value = 7809.727
$ws.Range("K122").Value = 18417.3999
$ws.Range("L122").Value = 23429.181
$ws.Range("M122").Value = -15967.3999
$ws.Range("N122").Value = -28329.181
$ws.Range("H138").Value = 84492.336
$ws.Range("J138").Value = 109238.5
$ws.Range("L138").Value = 109238.5
$ws.Range("N138").Value = -119518.5
$ws.Range("H141").Value = 116380.4
$ws.Range("J141").Value = 116380.4
$ws.Range("L141").Value = 116380.4
$ws.Range("N141").Value = -126740.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25008864
$ws.Range("J122").Value = 50008340
$ws.Range("L122").Value = 150025020
$ws.Range("N122").Value = -150029920
$ws.Range("H126").Value = 2937.3333
$ws.Range("I126").Value = 2712.3076
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 8136.9228
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -5666.9228
$ws.Range("N126").Value = -18140
$ws.Range("H137").Value = 75830.336
$ws.Range("J137").Value = 75830.336
$ws.Range("L137").Value = 75830.336
$ws.Range("N137").Value = -86030.336
$ws.Range("H140").Value = 149187.2
$ws.Range("J140").Value = 149187.2
$ws.Range("L140").Value = 149187.2
$ws.Range("N140").Value = -159547.2
